$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data pulled by the scheduled GitHub Actions job.
# Price cells in column D are textual (inline strings in the source data),
# using "." as a grouping separator, so numeric-looking values are written
# with a leading apostrophe to force Excel to store them as text (matching
# the original "General"-formatted text cells) instead of auto-converting
# them to numbers.

$ws.Range("D2").Value = "62.802.86"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.462.19"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'574.48"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'146.77"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "2.462.90"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'28.97"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "2.909.41"
$ws.Range("D17").Value = "62.780.35"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "2.465.97"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "'7.95"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'326.10"
$ws.Range("E21").Value = "  -1.26%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'2.22"
$ws.Range("E23").Value = "  +8.57%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'10.02"
$ws.Range("E25").Value = "  +18.01%  "
$ws.Range("D26").Value = "'65.54"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "0.0₃0986"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").Value = "2.582.89"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  -15.19%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "'7.98"
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "'152.09"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.369"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "'18.70"
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.38"
$ws.Range("E42").Value = "  -2.20%  "
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "0.0₆0309"
$ws.Range("E44").Value = "  -41.91%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'152.07"
$ws.Range("E46").Value = "  +4.66%  "
$ws.Range("E47").Value = "  +2.09%  "
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "'20.53"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("D50").Value = "'0.607"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("E51").Value = "  -1.19%  "
